$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 / J1 should carry the same formatting as the other
# header cells (e.g. H1), so copy H1's formats over before setting values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 2

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
